# Apply the "add 2022-Q3 data" edit:
#  1. Insert a new worksheet "2022-Q3" right after "总计" (pushing the other
#     quarter sheets one position later), by duplicating the existing
#     "2022-Q2" sheet so it keeps the same look & feel (borders, bold
#     header row, page margins, etc.), then overwriting its contents.
#  2. Populate "2022-Q3" with the new fund holdings data (9 rows).
#  3. Update the "总计" (summary) sheet to add a row for 2022-Q3 and shift
#     the previously existing rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: duplicate "2022-Q2" to create the new "2022-Q3" sheet in the
# right spot (right before "2022-Q2", i.e. right after "总计").
# ---------------------------------------------------------------------
$sheetQ2 = $wb.Worksheets.Item("2022-Q2")
$sheetQ2.Copy($sheetQ2)
$wsQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$wsQ3.Name = "2022-Q3"

# The template sheet only had 5 rows (1 header + 4 data); the new data has
# 8 data rows, so extend the existing formatting (borders/alignment on
# column A, "General" formatting for the rest) down to row 9 first.
$wsQ3.Range("A5:H5").Copy()
$wsQ3.Range("A6:H9").PasteSpecial(-4122) | Out-Null

# Force columns B:G to be stored as text for every data row, so values
# such as leading-zero fund codes or decimal-looking strings ("1.35",
# "012096", ...) are not silently reinterpreted as numbers by Excel.
$wsQ3.Range("B2:G9").NumberFormat = "@"

# ---------------------------------------------------------------------
# Step 2: fill in the "2022-Q3" fund holdings table.
# ---------------------------------------------------------------------
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$col = 2
foreach ($h in $headers) {
    $wsQ3.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# index, code, name, size, stockPos, posRatio, marketValue, posRank
$fundRows = @(
    @(0, "012096", "鑫元鑫动力混合A",       "2.04", "86.05", "4.09", "0.0834", 10),
    @(1, "014263", "鑫元长三角混合A",       "1.35", "78.27", "4.26", "0.0575", 5),
    @(2, "005262", "鑫元欣享灵活配置混合A", "0.86", "77.34", "3.72", "0.0320", 7),
    @(3, "012097", "鑫元鑫动力混合C",       "0.69", "86.05", "4.09", "0.0282", 10),
    @(4, "014264", "鑫元长三角混合C",       "0.34", "78.27", "4.26", "0.0145", 5),
    @(5, "012432", "国投瑞银安泰混合C",     "1.00", "32.06", "1.28", "0.0128", 8),
    @(6, "005263", "鑫元欣享灵活配置混合C", "0.22", "77.34", "3.72", "0.0082", 7),
    @(7, "012431", "国投瑞银安泰混合A",     "0.00", "32.06", "1.28", "0.0000", 8)
)

$r = 2
foreach ($item in $fundRows) {
    $wsQ3.Cells.Item($r, 1).Value = $item[0]
    $wsQ3.Cells.Item($r, 2).Value = $item[1]
    $wsQ3.Cells.Item($r, 3).Value = $item[2]
    $wsQ3.Cells.Item($r, 4).Value = $item[3]
    $wsQ3.Cells.Item($r, 5).Value = $item[4]
    $wsQ3.Cells.Item($r, 6).Value = $item[5]
    $wsQ3.Cells.Item($r, 7).Value = $item[6]
    $wsQ3.Cells.Item($r, 8).Value = $item[7]
    $r = $r + 1
}

# The very last row's "持有市值(亿元)" column (G) is actually a genuine
# number 0, not the text "0.0000" like the others.
$wsQ3.Range("G9").NumberFormat = "General"
$wsQ3.Range("G9").Value = 0

# ---------------------------------------------------------------------
# Step 3: update the "总计" summary sheet - add the 2022-Q3 row and shift
# the previously existing rows down by one.
# ---------------------------------------------------------------------
$wsTotal = $wb.Worksheets.Item("总计")

# index, quarter label, held count, held market value (亿元)
$totalData = @(
    @(0, "2022-Q3", 8, 0.24),
    @(1, "2022-Q2", 4, 0.07000000000000001),
    @(2, "2022-Q1", 4, 0.08),
    @(3, "2021-Q4", 5, 0.19),
    @(4, "2021-Q3", 7, 2.15)
)

$row = 2
foreach ($item in $totalData) {
    $wsTotal.Cells.Item($row, 1).Value = $item[0]
    $wsTotal.Cells.Item($row, 2).Value = $item[1]
    $wsTotal.Cells.Item($row, 3).Value = $item[2]
    $wsTotal.Cells.Item($row, 4).Value = $item[3]
    $row = $row + 1
}

# Column A keeps the same bordered/bold style on every data row (copy the
# style that was already used for the pre-existing A2 data cell).
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A2:A6").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# Restore the originally-selected tab (the last sheet, "2021-Q3") since
# duplicating a sheet makes the new copy the active/selected one.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2021-Q3").Activate()

Write-Host "2022-Q3 sheet inserted and summary sheet updated"
